$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" (second) and "C" (third) sub-rows within each year-group,
# columns A through E, to match the reordered source data.
# Year-group rows 3 (B) <-> 4 (C)
$A3_tmp = $ws.Range("A3").Value2
$B3_tmp = $ws.Range("B3").Value2
$C3_tmp = $ws.Range("C3").Value2
$D3_tmp = $ws.Range("D3").Value2
$E3_tmp = $ws.Range("E3").Value2
$A4_tmp = $ws.Range("A4").Value2
$B4_tmp = $ws.Range("B4").Value2
$C4_tmp = $ws.Range("C4").Value2
$D4_tmp = $ws.Range("D4").Value2
$E4_tmp = $ws.Range("E4").Value2
$ws.Range("A3").Value = $A4_tmp
$ws.Range("B3").Value = $B4_tmp
$ws.Range("C3").Value = $C4_tmp
$ws.Range("D3").Value = $D4_tmp
$ws.Range("E3").Value = $E4_tmp
$ws.Range("A4").Value = $A3_tmp
$ws.Range("B4").Value = $B3_tmp
$ws.Range("C4").Value = $C3_tmp
$ws.Range("D4").Value = $D3_tmp
$ws.Range("E4").Value = $E3_tmp

# Year-group rows 7 (B) <-> 8 (C)
$A7_tmp = $ws.Range("A7").Value2
$B7_tmp = $ws.Range("B7").Value2
$C7_tmp = $ws.Range("C7").Value2
$D7_tmp = $ws.Range("D7").Value2
$E7_tmp = $ws.Range("E7").Value2
$A8_tmp = $ws.Range("A8").Value2
$B8_tmp = $ws.Range("B8").Value2
$C8_tmp = $ws.Range("C8").Value2
$D8_tmp = $ws.Range("D8").Value2
$E8_tmp = $ws.Range("E8").Value2
$ws.Range("A7").Value = $A8_tmp
$ws.Range("B7").Value = $B8_tmp
$ws.Range("C7").Value = $C8_tmp
$ws.Range("D7").Value = $D8_tmp
$ws.Range("E7").Value = $E8_tmp
$ws.Range("A8").Value = $A7_tmp
$ws.Range("B8").Value = $B7_tmp
$ws.Range("C8").Value = $C7_tmp
$ws.Range("D8").Value = $D7_tmp
$ws.Range("E8").Value = $E7_tmp

# Year-group rows 11 (B) <-> 12 (C)
$A11_tmp = $ws.Range("A11").Value2
$B11_tmp = $ws.Range("B11").Value2
$C11_tmp = $ws.Range("C11").Value2
$D11_tmp = $ws.Range("D11").Value2
$E11_tmp = $ws.Range("E11").Value2
$A12_tmp = $ws.Range("A12").Value2
$B12_tmp = $ws.Range("B12").Value2
$C12_tmp = $ws.Range("C12").Value2
$D12_tmp = $ws.Range("D12").Value2
$E12_tmp = $ws.Range("E12").Value2
$ws.Range("A11").Value = $A12_tmp
$ws.Range("B11").Value = $B12_tmp
$ws.Range("C11").Value = $C12_tmp
$ws.Range("D11").Value = $D12_tmp
$ws.Range("E11").Value = $E12_tmp
$ws.Range("A12").Value = $A11_tmp
$ws.Range("B12").Value = $B11_tmp
$ws.Range("C12").Value = $C11_tmp
$ws.Range("D12").Value = $D11_tmp
$ws.Range("E12").Value = $E11_tmp

# Year-group rows 15 (B) <-> 16 (C)
$A15_tmp = $ws.Range("A15").Value2
$B15_tmp = $ws.Range("B15").Value2
$C15_tmp = $ws.Range("C15").Value2
$D15_tmp = $ws.Range("D15").Value2
$E15_tmp = $ws.Range("E15").Value2
$A16_tmp = $ws.Range("A16").Value2
$B16_tmp = $ws.Range("B16").Value2
$C16_tmp = $ws.Range("C16").Value2
$D16_tmp = $ws.Range("D16").Value2
$E16_tmp = $ws.Range("E16").Value2
$ws.Range("A15").Value = $A16_tmp
$ws.Range("B15").Value = $B16_tmp
$ws.Range("C15").Value = $C16_tmp
$ws.Range("D15").Value = $D16_tmp
$ws.Range("E15").Value = $E16_tmp
$ws.Range("A16").Value = $A15_tmp
$ws.Range("B16").Value = $B15_tmp
$ws.Range("C16").Value = $C15_tmp
$ws.Range("D16").Value = $D15_tmp
$ws.Range("E16").Value = $E15_tmp

# Year-group rows 19 (B) <-> 20 (C)
$A19_tmp = $ws.Range("A19").Value2
$B19_tmp = $ws.Range("B19").Value2
$C19_tmp = $ws.Range("C19").Value2
$D19_tmp = $ws.Range("D19").Value2
$E19_tmp = $ws.Range("E19").Value2
$A20_tmp = $ws.Range("A20").Value2
$B20_tmp = $ws.Range("B20").Value2
$C20_tmp = $ws.Range("C20").Value2
$D20_tmp = $ws.Range("D20").Value2
$E20_tmp = $ws.Range("E20").Value2
$ws.Range("A19").Value = $A20_tmp
$ws.Range("B19").Value = $B20_tmp
$ws.Range("C19").Value = $C20_tmp
$ws.Range("D19").Value = $D20_tmp
$ws.Range("E19").Value = $E20_tmp
$ws.Range("A20").Value = $A19_tmp
$ws.Range("B20").Value = $B19_tmp
$ws.Range("C20").Value = $C19_tmp
$ws.Range("D20").Value = $D19_tmp
$ws.Range("E20").Value = $E19_tmp

# Year-group rows 23 (B) <-> 24 (C)
$A23_tmp = $ws.Range("A23").Value2
$B23_tmp = $ws.Range("B23").Value2
$C23_tmp = $ws.Range("C23").Value2
$D23_tmp = $ws.Range("D23").Value2
$E23_tmp = $ws.Range("E23").Value2
$A24_tmp = $ws.Range("A24").Value2
$B24_tmp = $ws.Range("B24").Value2
$C24_tmp = $ws.Range("C24").Value2
$D24_tmp = $ws.Range("D24").Value2
$E24_tmp = $ws.Range("E24").Value2
$ws.Range("A23").Value = $A24_tmp
$ws.Range("B23").Value = $B24_tmp
$ws.Range("C23").Value = $C24_tmp
$ws.Range("D23").Value = $D24_tmp
$ws.Range("E23").Value = $E24_tmp
$ws.Range("A24").Value = $A23_tmp
$ws.Range("B24").Value = $B23_tmp
$ws.Range("C24").Value = $C23_tmp
$ws.Range("D24").Value = $D23_tmp
$ws.Range("E24").Value = $E23_tmp

# Year-group rows 27 (B) <-> 28 (C)
$A27_tmp = $ws.Range("A27").Value2
$B27_tmp = $ws.Range("B27").Value2
$C27_tmp = $ws.Range("C27").Value2
$D27_tmp = $ws.Range("D27").Value2
$E27_tmp = $ws.Range("E27").Value2
$A28_tmp = $ws.Range("A28").Value2
$B28_tmp = $ws.Range("B28").Value2
$C28_tmp = $ws.Range("C28").Value2
$D28_tmp = $ws.Range("D28").Value2
$E28_tmp = $ws.Range("E28").Value2
$ws.Range("A27").Value = $A28_tmp
$ws.Range("B27").Value = $B28_tmp
$ws.Range("C27").Value = $C28_tmp
$ws.Range("D27").Value = $D28_tmp
$ws.Range("E27").Value = $E28_tmp
$ws.Range("A28").Value = $A27_tmp
$ws.Range("B28").Value = $B27_tmp
$ws.Range("C28").Value = $C27_tmp
$ws.Range("D28").Value = $D27_tmp
$ws.Range("E28").Value = $E27_tmp

# Year-group rows 31 (B) <-> 32 (C)
$A31_tmp = $ws.Range("A31").Value2
$B31_tmp = $ws.Range("B31").Value2
$C31_tmp = $ws.Range("C31").Value2
$D31_tmp = $ws.Range("D31").Value2
$E31_tmp = $ws.Range("E31").Value2
$A32_tmp = $ws.Range("A32").Value2
$B32_tmp = $ws.Range("B32").Value2
$C32_tmp = $ws.Range("C32").Value2
$D32_tmp = $ws.Range("D32").Value2
$E32_tmp = $ws.Range("E32").Value2
$ws.Range("A31").Value = $A32_tmp
$ws.Range("B31").Value = $B32_tmp
$ws.Range("C31").Value = $C32_tmp
$ws.Range("D31").Value = $D32_tmp
$ws.Range("E31").Value = $E32_tmp
$ws.Range("A32").Value = $A31_tmp
$ws.Range("B32").Value = $B31_tmp
$ws.Range("C32").Value = $C31_tmp
$ws.Range("D32").Value = $D31_tmp
$ws.Range("E32").Value = $E31_tmp

# Year-group rows 35 (B) <-> 36 (C)
$A35_tmp = $ws.Range("A35").Value2
$B35_tmp = $ws.Range("B35").Value2
$C35_tmp = $ws.Range("C35").Value2
$D35_tmp = $ws.Range("D35").Value2
$E35_tmp = $ws.Range("E35").Value2
$A36_tmp = $ws.Range("A36").Value2
$B36_tmp = $ws.Range("B36").Value2
$C36_tmp = $ws.Range("C36").Value2
$D36_tmp = $ws.Range("D36").Value2
$E36_tmp = $ws.Range("E36").Value2
$ws.Range("A35").Value = $A36_tmp
$ws.Range("B35").Value = $B36_tmp
$ws.Range("C35").Value = $C36_tmp
$ws.Range("D35").Value = $D36_tmp
$ws.Range("E35").Value = $E36_tmp
$ws.Range("A36").Value = $A35_tmp
$ws.Range("B36").Value = $B35_tmp
$ws.Range("C36").Value = $C35_tmp
$ws.Range("D36").Value = $D35_tmp
$ws.Range("E36").Value = $E35_tmp

# Year-group rows 39 (B) <-> 40 (C)
$A39_tmp = $ws.Range("A39").Value2
$B39_tmp = $ws.Range("B39").Value2
$C39_tmp = $ws.Range("C39").Value2
$D39_tmp = $ws.Range("D39").Value2
$E39_tmp = $ws.Range("E39").Value2
$A40_tmp = $ws.Range("A40").Value2
$B40_tmp = $ws.Range("B40").Value2
$C40_tmp = $ws.Range("C40").Value2
$D40_tmp = $ws.Range("D40").Value2
$E40_tmp = $ws.Range("E40").Value2
$ws.Range("A39").Value = $A40_tmp
$ws.Range("B39").Value = $B40_tmp
$ws.Range("C39").Value = $C40_tmp
$ws.Range("D39").Value = $D40_tmp
$ws.Range("E39").Value = $E40_tmp
$ws.Range("A40").Value = $A39_tmp
$ws.Range("B40").Value = $B39_tmp
$ws.Range("C40").Value = $C39_tmp
$ws.Range("D40").Value = $D39_tmp
$ws.Range("E40").Value = $E39_tmp

# Year-group rows 43 (B) <-> 44 (C)
$A43_tmp = $ws.Range("A43").Value2
$B43_tmp = $ws.Range("B43").Value2
$C43_tmp = $ws.Range("C43").Value2
$D43_tmp = $ws.Range("D43").Value2
$E43_tmp = $ws.Range("E43").Value2
$A44_tmp = $ws.Range("A44").Value2
$B44_tmp = $ws.Range("B44").Value2
$C44_tmp = $ws.Range("C44").Value2
$D44_tmp = $ws.Range("D44").Value2
$E44_tmp = $ws.Range("E44").Value2
$ws.Range("A43").Value = $A44_tmp
$ws.Range("B43").Value = $B44_tmp
$ws.Range("C43").Value = $C44_tmp
$ws.Range("D43").Value = $D44_tmp
$ws.Range("E43").Value = $E44_tmp
$ws.Range("A44").Value = $A43_tmp
$ws.Range("B44").Value = $B43_tmp
$ws.Range("C44").Value = $C43_tmp
$ws.Range("D44").Value = $D43_tmp
$ws.Range("E44").Value = $E43_tmp

# Year-group rows 47 (B) <-> 48 (C)
$A47_tmp = $ws.Range("A47").Value2
$B47_tmp = $ws.Range("B47").Value2
$C47_tmp = $ws.Range("C47").Value2
$D47_tmp = $ws.Range("D47").Value2
$E47_tmp = $ws.Range("E47").Value2
$A48_tmp = $ws.Range("A48").Value2
$B48_tmp = $ws.Range("B48").Value2
$C48_tmp = $ws.Range("C48").Value2
$D48_tmp = $ws.Range("D48").Value2
$E48_tmp = $ws.Range("E48").Value2
$ws.Range("A47").Value = $A48_tmp
$ws.Range("B47").Value = $B48_tmp
$ws.Range("C47").Value = $C48_tmp
$ws.Range("D47").Value = $D48_tmp
$ws.Range("E47").Value = $E48_tmp
$ws.Range("A48").Value = $A47_tmp
$ws.Range("B48").Value = $B47_tmp
$ws.Range("C48").Value = $C47_tmp
$ws.Range("D48").Value = $D47_tmp
$ws.Range("E48").Value = $E47_tmp

# Year-group rows 51 (B) <-> 52 (C)
$A51_tmp = $ws.Range("A51").Value2
$B51_tmp = $ws.Range("B51").Value2
$C51_tmp = $ws.Range("C51").Value2
$D51_tmp = $ws.Range("D51").Value2
$E51_tmp = $ws.Range("E51").Value2
$A52_tmp = $ws.Range("A52").Value2
$B52_tmp = $ws.Range("B52").Value2
$C52_tmp = $ws.Range("C52").Value2
$D52_tmp = $ws.Range("D52").Value2
$E52_tmp = $ws.Range("E52").Value2
$ws.Range("A51").Value = $A52_tmp
$ws.Range("B51").Value = $B52_tmp
$ws.Range("C51").Value = $C52_tmp
$ws.Range("D51").Value = $D52_tmp
$ws.Range("E51").Value = $E52_tmp
$ws.Range("A52").Value = $A51_tmp
$ws.Range("B52").Value = $B51_tmp
$ws.Range("C52").Value = $C51_tmp
$ws.Range("D52").Value = $D51_tmp
$ws.Range("E52").Value = $E51_tmp

# Year-group rows 55 (B) <-> 56 (C)
$A55_tmp = $ws.Range("A55").Value2
$B55_tmp = $ws.Range("B55").Value2
$C55_tmp = $ws.Range("C55").Value2
$D55_tmp = $ws.Range("D55").Value2
$E55_tmp = $ws.Range("E55").Value2
$A56_tmp = $ws.Range("A56").Value2
$B56_tmp = $ws.Range("B56").Value2
$C56_tmp = $ws.Range("C56").Value2
$D56_tmp = $ws.Range("D56").Value2
$E56_tmp = $ws.Range("E56").Value2
$ws.Range("A55").Value = $A56_tmp
$ws.Range("B55").Value = $B56_tmp
$ws.Range("C55").Value = $C56_tmp
$ws.Range("D55").Value = $D56_tmp
$ws.Range("E55").Value = $E56_tmp
$ws.Range("A56").Value = $A55_tmp
$ws.Range("B56").Value = $B55_tmp
$ws.Range("C56").Value = $C55_tmp
$ws.Range("D56").Value = $D55_tmp
$ws.Range("E56").Value = $E55_tmp

# Year-group rows 59 (B) <-> 60 (C)
$A59_tmp = $ws.Range("A59").Value2
$B59_tmp = $ws.Range("B59").Value2
$C59_tmp = $ws.Range("C59").Value2
$D59_tmp = $ws.Range("D59").Value2
$E59_tmp = $ws.Range("E59").Value2
$A60_tmp = $ws.Range("A60").Value2
$B60_tmp = $ws.Range("B60").Value2
$C60_tmp = $ws.Range("C60").Value2
$D60_tmp = $ws.Range("D60").Value2
$E60_tmp = $ws.Range("E60").Value2
$ws.Range("A59").Value = $A60_tmp
$ws.Range("B59").Value = $B60_tmp
$ws.Range("C59").Value = $C60_tmp
$ws.Range("D59").Value = $D60_tmp
$ws.Range("E59").Value = $E60_tmp
$ws.Range("A60").Value = $A59_tmp
$ws.Range("B60").Value = $B59_tmp
$ws.Range("C60").Value = $C59_tmp
$ws.Range("D60").Value = $D59_tmp
$ws.Range("E60").Value = $E59_tmp

# Year-group rows 63 (B) <-> 64 (C)
$A63_tmp = $ws.Range("A63").Value2
$B63_tmp = $ws.Range("B63").Value2
$C63_tmp = $ws.Range("C63").Value2
$D63_tmp = $ws.Range("D63").Value2
$E63_tmp = $ws.Range("E63").Value2
$A64_tmp = $ws.Range("A64").Value2
$B64_tmp = $ws.Range("B64").Value2
$C64_tmp = $ws.Range("C64").Value2
$D64_tmp = $ws.Range("D64").Value2
$E64_tmp = $ws.Range("E64").Value2
$ws.Range("A63").Value = $A64_tmp
$ws.Range("B63").Value = $B64_tmp
$ws.Range("C63").Value = $C64_tmp
$ws.Range("D63").Value = $D64_tmp
$ws.Range("E63").Value = $E64_tmp
$ws.Range("A64").Value = $A63_tmp
$ws.Range("B64").Value = $B63_tmp
$ws.Range("C64").Value = $C63_tmp
$ws.Range("D64").Value = $D63_tmp
$ws.Range("E64").Value = $E63_tmp

# Remove the obsolete "gasoline production-sales ratio" (F) and
# "gasoline sales volume" (G) columns entirely.
$ws.Range("F1:G1").EntireColumn.Delete()

